$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1. Delete paragraphs that are removed outright (work from the
#    bottom of the document upward so earlier paragraph indices
#    stay valid).
# ---------------------------------------------------------------

# Lone empty paragraph that used to separate the "Replace with: \'"
# paragraph from the duplicate "Find what: Slide ..." block (this
# duplicate block is being moved further up, right after the
# bookmark paragraph).
$d.Range($d.Paragraphs(18).Range.Start, $d.Paragraphs(18).Range.End).Delete()

# The whole "Copy and paste entire document into Notepad..." /
# "Find what: ^13" / "Options: Use Wildcards" / "Replace with:  " /
# "[Replace line breaks with spaces]" / "Find what: '" block
# (paragraphs 8-16), which is no longer part of the process.
$d.Range($d.Paragraphs(8).Range.Start, $d.Paragraphs(16).Range.End).Delete()

# The two "Add "Slide 0" ..." / "Add "There are no speaker notes...""
# paragraphs (and their blank-line separators), paragraphs 3-6.
$d.Range($d.Paragraphs(3).Range.Start, $d.Paragraphs(6).Range.End).Delete()

# ---------------------------------------------------------------
# After the deletions above the document reads:
#   1: Find and replace TTS
#   2: (empty)
#   3: Convert superscripts and E's to "to the" for speaking exponents
#   4: (bookmark _GoBack) Replace with: \'
#   5: Find what: Slide [0-9]{1,}
#   6: Options: Use Wildcards
#   7: Replace with: ',^p'
#   8: (empty)
#   9: Delete ', from the front, add ' to the back
#   10: (empty)
#   11: Copy document into array
# ---------------------------------------------------------------

# 2. Append the parenthetical to the "Convert superscripts..." line.
$d.Paragraphs(3).Range.InsertXML('<w:p><w:r><w:t>Convert superscripts and E''s to "to the" for speaking exponents</w:t></w:r><w:r><w:t xml:space="preserve"> (search for Font: Superscript or [0-9</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>]E</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>, Use Wildcards)</w:t></w:r></w:p>')

# 3. Clear the run text on the bookmark paragraph, leaving only the
#    bookmark itself.
$d.Paragraphs(4).Range.InsertXML('<w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>')

# 4. Insert the (moved) "Find what: Slide [0-9]{1,}" / "Options: Use
#    Wildcards" / "Replace with: `,^p`" block right after the bookmark
#    paragraph.
$d.Paragraphs(4).Range.InsertParagraphAfter()
$d.Paragraphs(5).Range.InsertXML('<w:p><w:r><w:t xml:space="preserve">Find what: </w:t></w:r><w:r><w:t>Slide</w:t></w:r><w:r><w:t xml:space="preserve"> [0-9</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>]{</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>1,}</w:t></w:r></w:p>')

$d.Paragraphs(5).Range.InsertParagraphAfter()
$d.Paragraphs(6).Range.InsertXML('<w:p><w:r><w:t>Options: Use Wildcards</w:t></w:r></w:p>')

$d.Paragraphs(6).Range.InsertParagraphAfter()
$d.Paragraphs(7).Range.InsertXML('<w:p><w:r><w:t xml:space="preserve">Replace with: </w:t></w:r><w:r><w:t>`</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>,</w:t></w:r><w:r><w:t>^</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>p</w:t></w:r><w:r><w:t>`</w:t></w:r></w:p>')

# ---------------------------------------------------------------
# The document now reads:
#   1: Find and replace TTS
#   2: (empty)
#   3: Convert superscripts ... Use Wildcards)
#   4: (bookmark _GoBack) (empty)
#   5: Find what: Slide [0-9]{1,}
#   6: Options: Use Wildcards
#   7: Replace with: `,^p`
#   8: Find what: Slide [0-9]{1,}       <- old duplicate, to delete
#   9: Options: Use Wildcards            <- old duplicate, to delete
#   10: Replace with: ',^p'             <- old duplicate, to delete
#   11: (empty)
#   12: Delete ', from the front, add ' to the back
#   13: (empty)
#   14: Copy document into array
# ---------------------------------------------------------------

# 5. Remove the old duplicate "Find what: Slide ..." / "Options: Use
#    Wildcards" / "Replace with: ',^p'" block (paragraphs 8-10).
$d.Range($d.Paragraphs(8).Range.Start, $d.Paragraphs(10).Range.End).Delete()

# 6. Update the "Delete '... to the back" paragraph: swap the
#    straight quotes for backticks and split across two runs.
$d.Paragraphs(9).Range.InsertXML('<w:p><w:r><w:t>Delete `, from the front, add `</w:t></w:r><w:r><w:t xml:space="preserve"> to the back</w:t></w:r></w:p>')

Write-Host "Done. Paragraph count: $($d.Paragraphs.Count)"
